$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is inserted as row 9 (pushing the existing
# row 9..109 data down to rows 10..110, so the sheet grows from 109 to
# 110 rows). Insert a whole row first so every row below shifts down.
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C9").Value = 'Coquimbo'
$ws.Range("D9").Value = 44950
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 100112030
$ws.Range("G9").Value = 'Poroto granado'
$ws.Range("H9").Value = 'Sin especificar'
$ws.Range("I9").Value = 'Primera'
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 37000
$ws.Range("L9").Value = 38000
$ws.Range("M9").Value = 37500
$ws.Range("N9").Value = '$/malla 25 kilos'
$ws.Range("O9").Value = 'Provincia del Elquí'
$ws.Range("P9").Value = 1500
$ws.Range("Q9").Value = 25
$ws.Range("R9").Value = 'Hortaliza'

# Match the existing date column's number format on the new row (same
# style as the other rows in column D).
$ws.Range("D9").NumberFormat = $ws.Range("D10").NumberFormat
